$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "314.92"
Set-TextValue $ws.Range("E2") "1.98%"
Set-TextValue $ws.Range("D3") "39.22"
Set-TextValue $ws.Range("E3") "-1.45%"
Set-TextValue $ws.Range("D4") "5.146"
Set-TextValue $ws.Range("E4") "0.05%"
Set-TextValue $ws.Range("E5") "0.27%"
Set-TextValue $ws.Range("D6") "1.995"
Set-TextValue $ws.Range("E6") "2.50%"
Set-TextValue $ws.Range("D7") "8.343"
Set-TextValue $ws.Range("E7") "2.13%"
Set-TextValue $ws.Range("D8") "0.9364"
Set-TextValue $ws.Range("E8") "0.70%"
Set-TextValue $ws.Range("D9") "0.1301"
Set-TextValue $ws.Range("E9") "-9.49%"
Set-TextValue $ws.Range("D10") "0.1970"
Set-TextValue $ws.Range("E10") "2.48%"
Set-TextValue $ws.Range("D11") "0.08995"
Set-TextValue $ws.Range("E11") "-0.77%"
Set-TextValue $ws.Range("D12") "0.03523"
Set-TextValue $ws.Range("E12") "-0.04%"
Set-TextValue $ws.Range("D13") "0.09732"
Set-TextValue $ws.Range("E13") "-0.58%"
Set-TextValue $ws.Range("D14") "0.001405"
Set-TextValue $ws.Range("E14") "0.53%"
Set-TextValue $ws.Range("D15") "0.006132"
Set-TextValue $ws.Range("E15") "3.80%"
Set-TextValue $ws.Range("E16") "-7.56%"
Set-TextValue $ws.Range("E17") "3.70%"
Set-TextValue $ws.Range("D18") "3.137"
Set-TextValue $ws.Range("E18") "-6.25%"
Set-TextValue $ws.Range("D19") "0.3469"
Set-TextValue $ws.Range("E19") "1.13%"
Set-TextValue $ws.Range("D20") "0.1308"
Set-TextValue $ws.Range("E20") "-0.30%"
Set-TextValue $ws.Range("D21") "5.006"
Set-TextValue $ws.Range("E21") "6.69%"
Set-TextValue $ws.Range("D22") "0.2489"
Set-TextValue $ws.Range("E22") "2.61%"
Set-TextValue $ws.Range("D23") "0.04359"
Set-TextValue $ws.Range("E23") "-0.45%"
Set-TextValue $ws.Range("D24") "0.001239"
Set-TextValue $ws.Range("E24") "0.72%"
Set-TextValue $ws.Range("D25") "0.004761"
Set-TextValue $ws.Range("E25") "8.80%"
Set-TextValue $ws.Range("D26") "0.0003892"
Set-TextValue $ws.Range("E26") "199.30%"
Set-TextValue $ws.Range("E27") "-7.59%"
Set-TextValue $ws.Range("D39") "0.02228"
Set-TextValue $ws.Range("E39") "8.36%"
Set-TextValue $ws.Range("D40") "0.05200"
Set-TextValue $ws.Range("E40") "2.55%"
Set-TextValue $ws.Range("D41") "0.007745"
Set-TextValue $ws.Range("E41") "4.40%"
Set-TextValue $ws.Range("D42") "0.01030"
Set-TextValue $ws.Range("E42") "4.24%"
Set-TextValue $ws.Range("D43") "0.1398"
Set-TextValue $ws.Range("E43") "2.46%"
Set-TextValue $ws.Range("E44") "-1.38%"
Set-TextValue $ws.Range("D45") "0.008693"
Set-TextValue $ws.Range("E45") "-7.29%"
Set-TextValue $ws.Range("D46") "0.00006819"
Set-TextValue $ws.Range("E46") "6.63%"
Set-TextValue $ws.Range("E47") "-0.02%"
Set-TextValue $ws.Range("E48") "10.77%"
Set-TextValue $ws.Range("D49") "0.001691"
Set-TextValue $ws.Range("E49") "30.00%"
Set-TextValue $ws.Range("E50") "-0.02%"
Set-TextValue $ws.Range("E51") "-0.02%"
